$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# Update version string in A2
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# Update recommended citation string in A6
$oldCitation = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Huangling No. 2 Coal Mine, China, M0252, version ''' + $oldVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'
$newCitation = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Huangling No. 2 Coal Mine, China, M0252, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'
$aboutSheet.Range("A6").Value = $newCitation

# Update build_version column (S2:S10) on the data sheet
for ($row = 2; $row -le 10; $row++) {
    $dataSheet.Cells.Item($row, 19).Value = $newVersion
}
